# LTMC.MI.xlsx refresh — "Results from R script"
#
# The upstream R scraper re-ran and produced:
#   - a corrected close/adj_close pair for the existing last row (2024-06-11,
#     previously stamped with an intraday timestamp and stale close),
#   - the previously-missing 2024-06-12 daily bar, and
#   - the new 2024-06-13 daily bar.
#
# Net effect on the sheet: row 284 is rewritten in place and two brand-new
# rows (285, 286) are appended below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a genuine *text* value (shared-string backed,
# no number coercion, no stray NumberFormat left behind on the cell/style
# table) by round-tripping the text through a scratch cell's formula result
# and a values-only paste — mirrors how Range.PasteSpecial(xlPasteValues)
# behaves in real Excel COM automation.
function Set-TextValue {
    param($Row, $Col, [string]$Text)

    $scratch = $ws.Range("Z1")
    $escaped = $Text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Cells.Item($Row, $Col).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# Row 284 — rewritten: 2024-06-11 bar (date normalised, OHLC + adj_close
# refreshed). Ticker (H284) is already correct and untouched.
# ---------------------------------------------------------------------
$ws.Cells.Item(284, 1).Value = 45454.2916666667
$ws.Cells.Item(284, 2).Value = 106210
$ws.Cells.Item(284, 3).Value = 10.9499998092651
$ws.Cells.Item(284, 4).Value = 10.8400001525879
$ws.Cells.Item(284, 5).Value = 10.9499998092651
$ws.Cells.Item(284, 6).Value = 10.8699998855591
Set-TextValue 284 7 "10.8699998855591"

# ---------------------------------------------------------------------
# Row 285 — new: 2024-06-12 bar.
# ---------------------------------------------------------------------
$ws.Rows.Item(285).Insert()
$ws.Cells.Item(285, 1).Value = 45455.2916666667
$ws.Cells.Item(285, 2).Value = 180726
$ws.Cells.Item(285, 3).Value = 10.9300003051758
$ws.Cells.Item(285, 4).Value = 10.8500003814697
$ws.Cells.Item(285, 5).Value = 10.8900003433228
$ws.Cells.Item(285, 6).Value = 10.8900003433228
Set-TextValue 285 7 "10.8900003433228"
Set-TextValue 285 8 "LTMC.MI"

# ---------------------------------------------------------------------
# Row 286 — new: 2024-06-13 bar.
# ---------------------------------------------------------------------
$ws.Rows.Item(286).Insert()
$ws.Cells.Item(286, 1).Value = 45456.6495717593
$ws.Cells.Item(286, 2).Value = 131363
$ws.Cells.Item(286, 3).Value = 10.9499998092651
$ws.Cells.Item(286, 4).Value = 10.7299995422363
$ws.Cells.Item(286, 5).Value = 10.9499998092651
$ws.Cells.Item(286, 6).Value = 10.6800003051758
Set-TextValue 286 7 "10.6800003051758"
Set-TextValue 286 8 "LTMC.MI"
